# Normalize how <note> (Quote) is rendered: the lead paragraph used to be
# styled with the margin-note frame style (marginOuter + a direct
# w:framePr="around" override); it now uses a proper "Quote" paragraph
# style, and that style (plus its linked "Quote Char" character style) is
# defined in the stylesheet.

$d = $word.ActiveDocument

# --- 1. Re-point the lead paragraph at the "Quote" style and drop the
#        direct frame-wrap override it used to carry. -----------------
$xml = $d.WordOpenXML

$oldPPr = '<w:pPr><w:pStyle w:val="marginOuter"/><w:framePr w:wrap="around"/></w:pPr>'
$newPPr = '<w:pPr><w:pStyle w:val="Quote"/></w:pPr>'
if ($xml.Contains($oldPPr)) {
    $xml = $xml.Replace($oldPPr, $newPPr)
}

# --- 2. Define the "Quote" / "Quote Char" style pair (mirrors Word's
#        built-in Quote style) right after the document's last style. ---
$quoteStyles = '<w:style w:type="paragraph" w:styleId="Quote"><w:name w:val="Quote"/><w:basedOn w:val="Normal"/><w:next w:val="Normal"/><w:link w:val="QuoteChar"/><w:uiPriority w:val="29"/><w:semiHidden/><w:qFormat/><w:rsid w:val="005157A8"/><w:pPr><w:spacing w:before="120" w:after="120"/><w:ind w:left="113" w:right="113"/></w:pPr><w:rPr><w:i/><w:iCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:style><w:style w:type="character" w:customStyle="1" w:styleId="QuoteChar"><w:name w:val="Quote Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Quote"/><w:uiPriority w:val="29"/><w:semiHidden/><w:rsid w:val="005157A8"/><w:rPr><w:rFonts w:ascii="Baskerville" w:hAnsi="Baskerville"/><w:i/><w:iCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:style>'

$stylesClose = '</w:styles>'
if ($xml.Contains($stylesClose) -and -not $xml.Contains('w:styleId="Quote"')) {
    $lastClose = $xml.LastIndexOf($stylesClose)
    $xml = $xml.Substring(0, $lastClose) + $quoteStyles + $xml.Substring($lastClose)
}

$d.WordOpenXML = $xml

Write-Host "Quote style wired up; lead paragraph repointed."
